$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Discussion text for the two new meeting rows (written first so they land
# at the front of the shared-strings table, matching authoring order).
$ws.Range("C19").Value = "peer review "
$ws.Range("C21").Value = "Discussed about deliverable 3 and how we are gonna go about the demo on Thursday"

# Row 19 - new "peer review" meeting entry (attendance + duration)
$ws.Range("F19").Value = "1)abhinav"
$ws.Range("G19").Value = "2)Harsha"
$ws.Range("H19").Value = "3)Udhay"
$ws.Range("I19").Value = "4) Prakyath"
$ws.Range("K19").Value = "10:15-11:20"

# Row 21 - new deliverable 3 discussion entry (attendance + duration)
$ws.Range("F21").Value = "1)abhinav"
$ws.Range("G21").Value = "2)Harsha"
$ws.Range("H21").Value = "3)Udhay"
$ws.Range("I21").Value = "4)Prakyath"
$ws.Range("K21").Value = "10:15-11:20"

# Update selection to reflect the last-edited cell and view scroll
$ws.Range("K21").Select()
